$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (commessa 251218)
$ws.Range("A4").Value = 251218
$ws.Range("B4").Value = 'BIMEC 5'
$ws.Range("C4").Value = 23
$ws.Range("D4").Value = 96.90140845070422
$ws.Range("E4").Value = '2025-06-04 12:00:00'
$ws.Range("F4").Value = '2025-06-04 12:23:00'
$ws.Range("G4").Value = '2025-06-04 12:23:00'
$ws.Range("H4").Value = '2025-06-04 13:59:54'
$ws.Range("I4").Value = 6880
$ws.Range("J4").Value = 'bobina'
$ws.Range("K4").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R9'
$ws.Range("L4").Value = 6
$ws.Range("M4").Value = 76
$ws.Range("N4").Value = 39885
$ws.Range("O4").Value = 'X'
$ws.Range("P4").Value = 39885
$ws.Range("Q4").Value = '2025-05-09 00:00:00'
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 1

# Row 5 (commessa 251752)
$ws.Range("A5").Value = 251752
$ws.Range("B5").Value = 'BIMEC 5'
$ws.Range("C5").Value = 21
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = '2025-06-04 13:59:54'
$ws.Range("F5").Value = '2025-06-04 14:20:54'
$ws.Range("G5").Value = '2025-06-04 14:20:54'
$ws.Range("H5").Value = '2025-06-04 14:20:54'
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 'bobina'
$ws.Range("K5").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R6 ;R9'
$ws.Range("L5").Value = 3
$ws.Range("M5").Value = 76
$ws.Range("N5").Value = 39846
$ws.Range("O5").Value = 'X'
$ws.Range("P5").Value = 39846
$ws.Range("Q5").Value = '2025-05-20 00:00:00'
$ws.Range("R5").Value = -0.5978482003125
$ws.Range("S5").Value = 1

# Row 6 (commessa 251500)
$ws.Range("A6").Value = 251500
$ws.Range("B6").Value = 'BIMEC 5'
$ws.Range("C6").Value = 32
$ws.Range("D6").Value = 139.3802816901408
$ws.Range("E6").Value = '2025-06-04 14:20:54'
$ws.Range("F6").Value = '2025-06-04 14:52:54'
$ws.Range("G6").Value = '2025-06-04 14:52:54'
$ws.Range("H6").Value = '2025-06-05 09:12:16'
$ws.Range("I6").Value = 9896
$ws.Range("J6").Value = 'bobina'
$ws.Range("K6").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R6 ;R9'
$ws.Range("L6").Value = 4
$ws.Range("M6").Value = 70
$ws.Range("N6").Value = 39885
$ws.Range("O6").Value = 'X'
$ws.Range("P6").Value = 39885
$ws.Range("Q6").Value = '2025-05-26 00:00:00'
$ws.Range("R6").Value = -0.3835289514814815
$ws.Range("S6").Value = 2

# Row 7 (commessa 251565)
$ws.Range("A7").Value = 251565
$ws.Range("B7").Value = 'BIMEC 5'
$ws.Range("C7").Value = 19
$ws.Range("D7").Value = 176.7464788732394
$ws.Range("E7").Value = '2025-06-05 09:12:16'
$ws.Range("F7").Value = '2025-06-05 09:31:16'
$ws.Range("G7").Value = '2025-06-05 09:31:16'
$ws.Range("H7").Value = '2025-06-05 12:28:01'
$ws.Range("I7").Value = 12549
$ws.Range("J7").Value = 'bobina'
$ws.Range("K7").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9'
$ws.Range("L7").Value = 2
$ws.Range("M7").Value = 70
$ws.Range("N7").Value = 39885
$ws.Range("O7").Value = 'X'
$ws.Range("P7").Value = 39885
$ws.Range("Q7").Value = '2025-06-10 00:00:00'
$ws.Range("R7").Value = -0.5194640062615741
$ws.Range("S7").Value = 1

# Row 8 (commessa 251780)
$ws.Range("A8").Value = 251780
$ws.Range("B8").Value = 'BIMEC 5'
$ws.Range("C8").Value = 17
$ws.Range("D8").Value = 342.2394366197183
$ws.Range("E8").Value = '2025-06-05 12:28:01'
$ws.Range("F8").Value = '2025-06-05 12:45:01'
$ws.Range("G8").Value = '2025-06-05 12:45:01'
$ws.Range("H8").Value = '2025-06-06 10:27:16'
$ws.Range("I8").Value = 24299
$ws.Range("J8").Value = 'bobina'
$ws.Range("K8").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9'
$ws.Range("L8").Value = 3
$ws.Range("M8").Value = 70
$ws.Range("N8").Value = '39887 (esterno)'
$ws.Range("O8").Value = 'X'
$ws.Range("P8").Value = 39887
$ws.Range("Q8").Value = '2025-06-18 00:00:00'
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 7

# Row 9 (commessa 251180)
$ws.Range("A9").Value = 251180
$ws.Range("B9").Value = 'CASON'
$ws.Range("C9").Value = 32.5
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = '2025-06-05 07:00:00'
$ws.Range("F9").Value = '2025-06-05 07:32:30'
$ws.Range("G9").Value = '2025-06-05 07:32:30'
$ws.Range("H9").Value = '2025-06-05 07:32:30'
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 'bobina'
$ws.Range("K9").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9'
$ws.Range("L9").Value = 4
$ws.Range("M9").Value = 70
$ws.Range("N9").Value = '39887 (esterno)'
$ws.Range("O9").Value = 'X'
$ws.Range("P9").Value = 39887
$ws.Range("Q9").Value = '2025-05-20 00:00:00'
$ws.Range("R9").Value = -16.31423611111111
$ws.Range("S9").Value = 7

# Row 10 (commessa 251070)
$ws.Range("A10").Value = 251070
$ws.Range("B10").Value = 'CASON'
$ws.Range("C10").Value = 34.5
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = '2025-06-05 07:32:30'
$ws.Range("F10").Value = '2025-06-05 08:07:00'
$ws.Range("G10").Value = '2025-06-05 08:07:00'
$ws.Range("H10").Value = '2025-06-05 08:07:00'
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 'bobina'
$ws.Range("K10").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R6 ;R9'
$ws.Range("L10").Value = 6
$ws.Range("M10").Value = 70
$ws.Range("N10").Value = 39885
$ws.Range("O10").Value = 'X'
$ws.Range("P10").Value = 39885
$ws.Range("Q10").Value = '2025-03-28 00:00:00'
$ws.Range("R10").Value = -0.3381944444444445
$ws.Range("S10").Value = 2

# Row 11 (commessa 251773)
$ws.Range("A11").Value = 251773
$ws.Range("B11").Value = 'CASON'
$ws.Range("C11").Value = 32.5
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = '2025-06-05 08:07:00'
$ws.Range("F11").Value = '2025-06-05 08:39:30'
$ws.Range("G11").Value = '2025-06-05 08:39:30'
$ws.Range("H11").Value = '2025-06-05 08:39:30'
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 'bobina'
$ws.Range("K11").Value = 'CASON ;R6'
$ws.Range("L11").Value = 7
$ws.Range("M11").Value = 70
$ws.Range("N11").Value = 39874
$ws.Range("O11").Value = 'X'
$ws.Range("P11").Value = 39874
$ws.Range("Q11").Value = '2025-05-25 00:00:00'
$ws.Range("R11").Value = -0.3607638888888889
$ws.Range("S11").Value = 1

# Row 12 (commessa 251895)
$ws.Range("A12").Value = 251895
$ws.Range("B12").Value = 'CASON'
$ws.Range("C12").Value = 36.5
$ws.Range("D12").Value = 321.7090909090909
$ws.Range("E12").Value = '2025-06-05 08:39:30'
$ws.Range("F12").Value = '2025-06-05 09:16:00'
$ws.Range("G12").Value = '2025-06-05 09:16:00'
$ws.Range("H12").Value = '2025-06-05 14:37:42'
$ws.Range("I12").Value = 17694
$ws.Range("J12").Value = 'bobina'
$ws.Range("K12").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R6 ;R9'
$ws.Range("L12").Value = 10
$ws.Range("M12").Value = 70
$ws.Range("N12").Value = '39891 (esterno)'
$ws.Range("O12").Value = 'X'
$ws.Range("P12").Value = 39891
$ws.Range("Q12").Value = '2025-05-26 00:00:00'
$ws.Range("R12").Value = -10.60952020202546
$ws.Range("S12").Value = 4

# Row 13 (commessa 252112)
$ws.Range("A13").Value = 252112
$ws.Range("B13").Value = 'R10'
$ws.Range("C13").Value = 20
$ws.Range("D13").Value = 204.1475409836065
$ws.Range("E13").Value = '2025-06-04 07:00:00'
$ws.Range("F13").Value = '2025-06-04 07:20:00'
$ws.Range("G13").Value = '2025-06-04 07:20:00'
$ws.Range("H13").Value = '2025-06-04 10:44:08'
$ws.Range("I13").Value = 12453
$ws.Range("J13").Value = 'bobina'
$ws.Range("K13").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9'
$ws.Range("L13").Value = 2
$ws.Range("M13").Value = 70
$ws.Range("N13").Value = 39885
$ws.Range("O13").Value = 'X'
$ws.Range("P13").Value = 39885
$ws.Range("Q13").Value = '2025-06-06 00:00:00'
$ws.Range("R13").Value = 0
$ws.Range("S13").Value = 1

# Row 14 (commessa 252282)
$ws.Range("A14").Value = 252282
$ws.Range("B14").Value = 'R10'
$ws.Range("C14").Value = 35
$ws.Range("D14").Value = 44.88524590163934
$ws.Range("E14").Value = '2025-06-04 10:44:08'
$ws.Range("F14").Value = '2025-06-04 11:19:08'
$ws.Range("G14").Value = '2025-06-04 11:19:08'
$ws.Range("H14").Value = '2025-06-04 12:04:01'
$ws.Range("I14").Value = 2738
$ws.Range("J14").Value = 'bobina'
$ws.Range("K14").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9'
$ws.Range("L14").Value = 5
$ws.Range("M14").Value = 70
$ws.Range("N14").Value = 39885
$ws.Range("O14").Value = 'X'
$ws.Range("P14").Value = 39885
$ws.Range("Q14").Value = '2025-06-09 00:00:00'
$ws.Range("R14").Value = 0
$ws.Range("S14").Value = 1

# Row 15 (commessa 251984)
$ws.Range("A15").Value = 251984
$ws.Range("B15").Value = 'R10'
$ws.Range("C15").Value = 30
$ws.Range("D15").Value = 338.327868852459
$ws.Range("E15").Value = '2025-06-04 12:04:01'
$ws.Range("F15").Value = '2025-06-04 12:34:01'
$ws.Range("G15").Value = '2025-06-04 12:34:01'
$ws.Range("H15").Value = '2025-06-05 10:12:21'
$ws.Range("I15").Value = 20638
$ws.Range("J15").Value = 'bobina'
$ws.Range("K15").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9'
$ws.Range("L15").Value = 3
$ws.Range("M15").Value = 70
$ws.Range("N15").Value = 39874
$ws.Range("O15").Value = 'X'
$ws.Range("P15").Value = 39874
$ws.Range("Q15").Value = '2025-06-10 00:00:00'
$ws.Range("R15").Value = -0.4252504553703704
$ws.Range("S15").Value = 1

# Row 16 (commessa 251362)
$ws.Range("A16").Value = 251362
$ws.Range("B16").Value = 'R12'
$ws.Range("C16").Value = 17
$ws.Range("D16").Value = 35.28169014084507
$ws.Range("E16").Value = '2025-06-05 07:00:00'
$ws.Range("F16").Value = '2025-06-05 07:17:00'
$ws.Range("G16").Value = '2025-06-05 07:17:00'
$ws.Range("H16").Value = '2025-06-05 07:52:16'
$ws.Range("I16").Value = 2505
$ws.Range("J16").Value = 'bobina'
$ws.Range("K16").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R9'
$ws.Range("L16").Value = 3
$ws.Range("M16").Value = 76
$ws.Range("N16").Value = 39874
$ws.Range("O16").Value = 'X'
$ws.Range("P16").Value = 39874
$ws.Range("Q16").Value = '2025-04-24 00:00:00'
$ws.Range("R16").Value = -0.3279733959259259
$ws.Range("S16").Value = 7

# Row 17 (commessa 251631)
$ws.Range("A17").Value = 251631
$ws.Range("B17").Value = 'R12'
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 146.056338028169
$ws.Range("E17").Value = '2025-06-05 07:52:16'
$ws.Range("F17").Value = '2025-06-05 08:07:16'
$ws.Range("G17").Value = '2025-06-05 08:07:16'
$ws.Range("H17").Value = '2025-06-05 10:33:20'
$ws.Range("I17").Value = 10370
$ws.Range("J17").Value = 'bobina'
$ws.Range("K17").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R6 ;R9'
$ws.Range("L17").Value = 3
$ws.Range("M17").Value = 76
$ws.Range("N17").Value = 39885
$ws.Range("O17").Value = 'X'
$ws.Range("P17").Value = 39885
$ws.Range("Q17").Value = '2025-05-09 00:00:00'
$ws.Range("R17").Value = -0.4398180751157408
$ws.Range("S17").Value = 2

# Row 18 (commessa 251237)
$ws.Range("A18").Value = 251237
$ws.Range("B18").Value = 'R12'
$ws.Range("C18").Value = 36
$ws.Range("D18").Value = 565.3521126760563
$ws.Range("E18").Value = '2025-06-05 10:33:20'
$ws.Range("F18").Value = '2025-06-05 11:09:20'
$ws.Range("G18").Value = '2025-06-05 11:09:20'
$ws.Range("H18").Value = '2025-06-06 12:34:41'
$ws.Range("I18").Value = 40140
$ws.Range("J18").Value = 'bobina'
$ws.Range("K18").Value = 'R12 ;R9'
$ws.Range("L18").Value = 6
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = 39885
$ws.Range("O18").Value = 'X'
$ws.Range("P18").Value = 39885
$ws.Range("Q18").Value = '2025-05-05 00:00:00'
$ws.Range("R18").Value = -1.524090375590278
$ws.Range("S18").Value = 1

# Row 19 (commessa 252084)
$ws.Range("A19").Value = 252084
$ws.Range("B19").Value = 'R3'
$ws.Range("C19").Value = 40
$ws.Range("D19").Value = 797.9795918367347
$ws.Range("E19").Value = '2025-06-04 07:00:00'
$ws.Range("F19").Value = '2025-06-04 07:40:00'
$ws.Range("G19").Value = '2025-06-04 07:40:00'
$ws.Range("H19").Value = '2025-06-05 12:57:58'
$ws.Range("I19").Value = 39101
$ws.Range("J19").Value = 'bobina'
$ws.Range("K19").Value = 'BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9'
$ws.Range("L19").Value = 2
$ws.Range("M19").Value = 70
$ws.Range("N19").Value = 39885
$ws.Range("O19").Value = 'X'
$ws.Range("P19").Value = 39885
$ws.Range("Q19").Value = '2025-06-30 00:00:00'
$ws.Range("R19").Value = -0.5402636054398148
$ws.Range("S19").Value = 7
